$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 47
$ws.Range("B2").Value = 121
$ws.Range("B3").Value = 142
$ws.Range("B4").Value = 169
$ws.Range("B5").Value = 199
$ws.Range("B6").Value = 207
